$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2674638667971811
$ws.Range("C2").Value = 1.761278749491586
$ws.Range("D2").Value = 11.53428961468776
$ws.Range("E2").Value = 3.39621695636303
$ws.Range("F2").Value = 3.465342416133991
$ws.Range("G2").Value = 22
